$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 89561920
$ws.Range("B2").Value = 77506
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 508180.1397883456
$ws.Range("R2").Value = 7157869.144381812
$ws.Range("A3").Value = 89561903
$ws.Range("B3").Value = 73693
$ws.Range("E3").Value = 6440
$ws.Range("F3").Value = "Vitgrynig nållav"
$ws.Range("G3").Value = "Chaenotheca subroscida"
$ws.Range("H3").Value = "(Eitner) Zahlbr."
$ws.Range("Q3").Value = 508181.8907144414
$ws.Range("R3").Value = 7157858.822880358
$ws.Range("A4").Value = 89561919
$ws.Range("B4").Value = 89392
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 508181.0314439956
$ws.Range("R4").Value = 7157857.960051162
$ws.Range("A5").Value = 110694996
$ws.Range("B5").Value = 77515
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("P5").Value = "Lill-bergvattnet, Jmt"
$ws.Range("Q5").Value = 508120.4989547321
$ws.Range("R5").Value = 7157940.836871861
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-07-07"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-07-07"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AW5").Value = "Ulrika Westling"
$ws.Range("AX5").Value = "Ulrika Westling"
$ws.Range("AY5").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A6").Value = 110694995
$ws.Range("B6").Value = 81248
$ws.Range("E6").Value = 1312
$ws.Range("F6").Value = "Gammelgransskål"
$ws.Range("G6").Value = "Pseudographis pinicola"
$ws.Range("H6").Value = "(Nyl.) Rehm"
$ws.Range("P6").Value = "Lill-Bergvattnet, Jmt"
$ws.Range("Q6").Value = 508120.5046955775
$ws.Range("R6").Value = 7157938.685647392
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-07-07"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-07-07"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AW6").Value = "Ulrika Westling"
$ws.Range("AX6").Value = "Ulrika Westling"
$ws.Range("AY6").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A7").Value = 110694985
$ws.Range("B7").Value = 56543
$ws.Range("E7").Value = 103021
$ws.Range("F7").Value = "Talltita"
$ws.Range("G7").Value = "Poecile montanus"
$ws.Range("H7").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("P7").Value = "Lill-Bergvattnet, Jmt"
$ws.Range("Q7").Value = 508117.4027607946
$ws.Range("R7").Value = 7157809.60020776
$ws.Range("S7").Value = 25
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-07-07"
$ws.Range("Y7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-07-07"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AC7").Value = "födosökande och tystlåten för en ggs skull"
$ws.Range("AW7").Value = "Ulrika Westling"
$ws.Range("AX7").Value = "Ulrika Westling"
$ws.Range("AY7").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A8").Value = 110694994
$ws.Range("B8").Value = 89423
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma"
$ws.Range("H8").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P8").Value = "Lill-bergvattnet, Jmt"
$ws.Range("Q8").Value = 508122.4460188448
$ws.Range("R8").Value = 7157856.942592324
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-07-07"
$ws.Range("Y8").Style = "Normal"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-07-07"
$ws.Range("AA8").Style = "Normal"
$ws.Range("AW8").Value = "Ulrika Westling"
$ws.Range("AX8").Value = "Ulrika Westling"
$ws.Range("AY8").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A9").Value = 110694992
$ws.Range("B9").Value = 73696
$ws.Range("P9").Value = "Lill-bergvattnet, Jmt"
$ws.Range("Q9").Value = 508120.7791116443
$ws.Range("R9").Value = 7157835.855348294
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-07-07"
$ws.Range("Y9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-07-07"
$ws.Range("AA9").Style = "Normal"
$ws.Range("AW9").Value = "Ulrika Westling"
$ws.Range("AX9").Value = "Ulrika Westling"
$ws.Range("AY9").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A10").Value = 89561904
$ws.Range("B10").Value = 89356
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 5447
$ws.Range("F10").Value = "Vedticka"
$ws.Range("G10").Value = "Fuscoporia viticola"
$ws.Range("H10").Value = "(Schwein.) Murrill"
$ws.Range("Q10").Value = 508405.8014842027
$ws.Range("R10").Value = 7157896.004994209
$ws.Range("A11").Value = 89561918
$ws.Range("B11").Value = 77506
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 508584.965041342
$ws.Range("R11").Value = 7157911.994302315
$ws.Range("AC11").Value = ""
$ws.Range("A12").Value = 89561888
$ws.Range("B12").Value = 81236
$ws.Range("E12").Value = 1312
$ws.Range("F12").Value = "Gammelgransskål"
$ws.Range("G12").Value = "Pseudographis pinicola"
$ws.Range("H12").Value = "(Nyl.) Rehm"
$ws.Range("P12").Value = "Öster Nåsjön, Jmt"
$ws.Range("Q12").Value = 508419.1011138954
$ws.Range("R12").Value = 7157915.833398769
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2020-07-09"
$ws.Range("Y12").Style = "Normal"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2020-07-09"
$ws.Range("AA12").Style = "Normal"
$ws.Range("AW12").Value = "Erland Lindblad"
$ws.Range("AX12").Value = "Via Erland Lindblad"
$ws.Range("AY12").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
$ws.Range("A13").Value = 89561889
$ws.Range("B13").Value = 77506
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("P13").Value = "Öster Nåsjön, Jmt"
$ws.Range("Q13").Value = 508399.8642952705
$ws.Range("R13").Value = 7157861.998314789
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2020-07-09"
$ws.Range("Y13").Style = "Normal"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2020-07-09"
$ws.Range("AA13").Style = "Normal"
$ws.Range("AW13").Value = "Erland Lindblad"
$ws.Range("AX13").Value = "Via Erland Lindblad"
$ws.Range("AY13").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
$ws.Range("A14").Value = 89561892
$ws.Range("B14").Value = 73693
$ws.Range("E14").Value = 6440
$ws.Range("F14").Value = "Vitgrynig nållav"
$ws.Range("G14").Value = "Chaenotheca subroscida"
$ws.Range("H14").Value = "(Eitner) Zahlbr."
$ws.Range("P14").Value = "Öster Nåsjön, Jmt"
$ws.Range("Q14").Value = 508399.8642952705
$ws.Range("R14").Value = 7157861.998314789
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2020-07-09"
$ws.Range("Y14").Style = "Normal"
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2020-07-09"
$ws.Range("AA14").Style = "Normal"
$ws.Range("AW14").Value = "Erland Lindblad"
$ws.Range("AX14").Value = "Via Erland Lindblad"
$ws.Range("AY14").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
$ws.Range("A15").Value = 89561876
$ws.Range("B15").Value = 77506
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("P15").Value = "Öster Nåsjön, Jmt"
$ws.Range("Q15").Value = 508556.0323171288
$ws.Range("R15").Value = 7157936.867342628
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2020-07-09"
$ws.Range("Y15").Style = "Normal"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2020-07-09"
$ws.Range("AA15").Style = "Normal"
$ws.Range("AW15").Value = "Erland Lindblad"
$ws.Range("AX15").Value = "Via Erland Lindblad"
$ws.Range("AY15").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
$ws.Range("A16").Value = 89561908
$ws.Range("B16").Value = 56395
$ws.Range("C16").Value = "Behöver inte valideras"
$ws.Range("E16").Value = 100109
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("P16").Value = "Öster Nåsjön, Jmt"
$ws.Range("Q16").Value = 508496.1585720535
$ws.Range("R16").Value = 7157934.978585394
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2020-07-09"
$ws.Range("Y16").Style = "Normal"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2020-07-09"
$ws.Range("AA16").Style = "Normal"
$ws.Range("AC16").Value = "Ringhack"
$ws.Range("AW16").Value = "Erland Lindblad"
$ws.Range("AX16").Value = "Via Erland Lindblad"
$ws.Range("AY16").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
$ws.Range("A17").Value = 110694989
$ws.Range("B17").Value = 89405
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 508301.697820781
$ws.Range("R17").Value = 7157841.936939664
$ws.Range("A18").Value = 110694988
$ws.Range("B18").Value = 89405
$ws.Range("E18").Value = 1202
$ws.Range("F18").Value = "Ullticka"
$ws.Range("G18").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H18").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q18").Value = 508407.9625505499
$ws.Range("R18").Value = 7157893.429437263
$ws.Range("A19").Value = 110694986
$ws.Range("B19").Value = 85715
$ws.Range("E19").Value = 510
$ws.Range("F19").Value = "Doftskinn"
$ws.Range("G19").Value = "Cystostereum murrayi"
$ws.Range("H19").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("P19").Value = "Lill-bergvattnet, Jmt"
$ws.Range("Q19").Value = 508339.86753408
$ws.Range("R19").Value = 7157904.428650577
$ws.Range("S19").Value = 10
$ws.Range("AC19").Value = ""
$ws.Range("A20").Value = 110694991
$ws.Range("B20").Value = 89419
$ws.Range("E20").Value = 1204
$ws.Range("F20").Value = "Gränsticka"
$ws.Range("G20").Value = "Phellopilus nigrolimitatus"
$ws.Range("H20").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q20").Value = 508317.2746237174
$ws.Range("R20").Value = 7157817.024181078
$ws.Range("A21").Value = 110694993
$ws.Range("Q21").Value = 508476.5727480573
$ws.Range("R21").Value = 7157852.31504698
$ws.Range("A22").Value = 110694987
$ws.Range("B22").Value = 85715
$ws.Range("E22").Value = 510
$ws.Range("F22").Value = "Doftskinn"
$ws.Range("G22").Value = "Cystostereum murrayi"
$ws.Range("H22").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q22").Value = 508440.3121548754
$ws.Range("R22").Value = 7157878.890296321
$ws.Range("A23").Value = 110694990
$ws.Range("B23").Value = 89419
$ws.Range("E23").Value = 1204
$ws.Range("F23").Value = "Gränsticka"
$ws.Range("G23").Value = "Phellopilus nigrolimitatus"
$ws.Range("H23").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q23").Value = 508432.1223925821
$ws.Range("R23").Value = 7157880.588618397
